$d = $word.ActiveDocument
$d.Content.Find.Execute("IntelliJ IDEA Ultimate 2018.3.3", $false, $false, $false, $false, $false, $true, 1, $false, "IntelliJ IDEA Ultimate 2018.3.5", 2)
